function Set-TextValue {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple in-place updates ---
Set-TextValue $ws.Cells.Item(2,4) "96.283.63"
Set-TextValue $ws.Cells.Item(2,5) "  +0.58%  "

Set-TextValue $ws.Cells.Item(3,4) "3.584.33"
Set-TextValue $ws.Cells.Item(3,5) "  -0.66%  "

Set-TextValue $ws.Cells.Item(4,5) "  -0.09%  "

Set-TextValue $ws.Cells.Item(5,4) "240.85"
Set-TextValue $ws.Cells.Item(5,5) "  +0.37%  "

Set-TextValue $ws.Cells.Item(6,4) "656.17"
Set-TextValue $ws.Cells.Item(6,5) "  +0.42%  "

Set-TextValue $ws.Cells.Item(7,5) "  +5.65%  "

Set-TextValue $ws.Cells.Item(8,5) "  -0.61%  "

Set-TextValue $ws.Cells.Item(9,5) "  +0.02%  "

Set-TextValue $ws.Cells.Item(10,4) "1.04"
Set-TextValue $ws.Cells.Item(10,5) "  +3.27%  "

Set-TextValue $ws.Cells.Item(11,4) "3.584.12"
Set-TextValue $ws.Cells.Item(11,5) "  -0.62%  "

Set-TextValue $ws.Cells.Item(12,4) "43.42"
Set-TextValue $ws.Cells.Item(12,5) "  +0.11%  "

Set-TextValue $ws.Cells.Item(13,5) "  +0.56%  "

Set-TextValue $ws.Cells.Item(14,4) "6.40"
Set-TextValue $ws.Cells.Item(14,5) "  +0.54%  "

Set-TextValue $ws.Cells.Item(15,4) "4.249.29"
Set-TextValue $ws.Cells.Item(15,5) "  -1.18%  "

Set-TextValue $ws.Cells.Item(16,4) "96.173.61"
Set-TextValue $ws.Cells.Item(16,5) "  +0.55%  "

Set-TextValue $ws.Cells.Item(17,5) "  +0.35%  "

Set-TextValue $ws.Cells.Item(18,4) "3.570.60"
Set-TextValue $ws.Cells.Item(18,5) "  -1.29%  "

Set-TextValue $ws.Cells.Item(19,5) "  -2.52%  "

Set-TextValue $ws.Cells.Item(20,4) "12.60"
Set-TextValue $ws.Cells.Item(20,5) "  +0.42%  "

Set-TextValue $ws.Cells.Item(21,5) "  -2.16%  "

Set-TextValue $ws.Cells.Item(22,4) "0.493"
Set-TextValue $ws.Cells.Item(22,5) "  +1.23%  "

Set-TextValue $ws.Cells.Item(25,4) "0.0000200"
Set-TextValue $ws.Cells.Item(25,5) "  +1.32%  "

Set-TextValue $ws.Cells.Item(26,4) "6.85"
Set-TextValue $ws.Cells.Item(26,5) "  +2.61%  "

Set-TextValue $ws.Cells.Item(27,4) "96.52"
Set-TextValue $ws.Cells.Item(27,5) "  -0.69%  "

Set-TextValue $ws.Cells.Item(28,4) "12.82"
Set-TextValue $ws.Cells.Item(28,5) "  -0.45%  "

Set-TextValue $ws.Cells.Item(29,4) "3.776.08"
Set-TextValue $ws.Cells.Item(29,5) "  -0.59%  "

Set-TextValue $ws.Cells.Item(30,4) "2.99"
Set-TextValue $ws.Cells.Item(30,5) "  -7.38%  "

Set-TextValue $ws.Cells.Item(31,4) "0.150"
Set-TextValue $ws.Cells.Item(31,5) "  +7.05%  "

Set-TextValue $ws.Cells.Item(32,5) "  +0.92%  "

Set-TextValue $ws.Cells.Item(33,5) "  +0.08%  "

Set-TextValue $ws.Cells.Item(34,4) "0.184"
Set-TextValue $ws.Cells.Item(34,5) "  +3.22%  "

Set-TextValue $ws.Cells.Item(35,4) "0.997"
Set-TextValue $ws.Cells.Item(35,5) "  -0.50%  "

Set-TextValue $ws.Cells.Item(36,4) "31.59"
Set-TextValue $ws.Cells.Item(36,5) "  -1.11%  "

Set-TextValue $ws.Cells.Item(39,4) "8.69"
Set-TextValue $ws.Cells.Item(39,5) "  +3.67%  "

Set-TextValue $ws.Cells.Item(40,4) "1.64"
Set-TextValue $ws.Cells.Item(40,5) "  +8.72%  "

Set-TextValue $ws.Cells.Item(41,5) "  +0.04%  "

Set-TextValue $ws.Cells.Item(42,5) "  -0.43%  "

Set-TextValue $ws.Cells.Item(43,4) "0.910"
Set-TextValue $ws.Cells.Item(43,5) "  -2.52%  "

Set-TextValue $ws.Cells.Item(44,4) "1.83"
Set-TextValue $ws.Cells.Item(44,5) "  +4.59%  "

Set-TextValue $ws.Cells.Item(45,5) "  -0.61%  "

Set-TextValue $ws.Cells.Item(46,5) "  +0.80%  "

Set-TextValue $ws.Cells.Item(49,4) "0.0418"
Set-TextValue $ws.Cells.Item(49,5) "  -0.51%  "

Set-TextValue $ws.Cells.Item(50,4) "3.58"
Set-TextValue $ws.Cells.Item(50,5) "  +3.17%  "

Set-TextValue $ws.Cells.Item(51,4) "53.41"
Set-TextValue $ws.Cells.Item(51,5) "  -2.09%  "

# --- Row swaps (content moved between adjacent rows) ---
# Row 23
Set-TextValue $ws.Cells.Item(23,2) "SuiNetwork"
Set-TextValue $ws.Cells.Item(23,3) "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Cells.Item(23,4) "3.45"
Set-TextValue $ws.Cells.Item(23,5) "  -2.23%  "

# Row 24
Set-TextValue $ws.Cells.Item(24,2) "BitcoinCash"
Set-TextValue $ws.Cells.Item(24,3) "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Cells.Item(24,4) "511.78"
Set-TextValue $ws.Cells.Item(24,5) "  -0.14%  "

# Row 37
Set-TextValue $ws.Cells.Item(37,2) "Bittensor"
Set-TextValue $ws.Cells.Item(37,3) "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Cells.Item(37,4) "614.65"
Set-TextValue $ws.Cells.Item(37,5) "  +7.85%  "

# Row 38
Set-TextValue $ws.Cells.Item(38,2) "PolygonEcosystemToken"
Set-TextValue $ws.Cells.Item(38,3) "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws.Cells.Item(38,4) "0.566"
Set-TextValue $ws.Cells.Item(38,5) "  +0.30%  "

# Row 47
Set-TextValue $ws.Cells.Item(47,2) "EnergySwap"
Set-TextValue $ws.Cells.Item(47,3) "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Cells.Item(47,4) "34.23"
Set-TextValue $ws.Cells.Item(47,5) "  -0.94%  "

# Row 48
Set-TextValue $ws.Cells.Item(48,2) "WhiteBITCoin"
Set-TextValue $ws.Cells.Item(48,3) "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Cells.Item(48,4) "23.51"
Set-TextValue $ws.Cells.Item(48,5) "  -1.20%  "
